$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.147.66"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "1.784.47"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("E6").Value = "  +1.19%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0687"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").Value = "2.042.94"
$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.819.48"
$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.622"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").Value = "34.055.39"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.37%  "

$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("E23").Value = "  +2.52%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.80%  "

$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("E31").Value = "  +2.14%  "

$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("D35").Value = "1.445.35"
$ws.Range("E35").Value = "  +4.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.654"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.93%  "

$ws.Range("E38").Value = "  +4.59%  "

$ws.Range("E39").Value = "  +1.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.98%  "

$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.921"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.74%  "

$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("D49").Value = "1.944.72"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("E51").Value = "  +0.04%  "
